$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.580.76'
$ws.Range("E2").Value = '  +1.54%  '

$ws.Range("D3").Value = '1.475.49'
$ws.Range("E3").Value = '  +2.18%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").Value = '0.9578'
$ws.Range("E5").Value = '  +3.36%  '

$ws.Range("D6").Value = '277.47'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3540'
$ws.Range("E7").Value = '  -2.81%  '

$ws.Range("D8").Value = '0.3075'
$ws.Range("E8").Value = '  -0.22%  '

$ws.Range("D9").Value = '1.083'
$ws.Range("E9").Value = '  +5.67%  '

$ws.Range("D10").Value = '39.34'
$ws.Range("E10").Value = '  -0.35%  '

$ws.Range("D11").Value = '0.06644'
$ws.Range("E11").Value = '  +1.77%  '

$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.56%  '

$ws.Range("D13").Value = '5.469'
$ws.Range("E13").Value = '  +2.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.10'
$ws.Range("E14").Value = '  +3.19%  '

$ws.Range("D15").Value = '6.172'
$ws.Range("E15").Value = '  +1.72%  '

$ws.Range("D16").Value = '0.9591'
$ws.Range("E16").Value = '  +1.40%  '

$ws.Range("D17").Value = '0.00001016'
$ws.Range("E17").Value = '  +0.31%  '

$ws.Range("D18").Value = '1.473.83'
$ws.Range("E18").Value = '  +2.42%  '

$ws.Range("D19").Value = '0.05981'
$ws.Range("E19").Value = '  +5.68%  '

$ws.Range("D20").Value = '68.91'
$ws.Range("E20").Value = '  -0.11%  '

$ws.Range("D21").Value = '5.486'
$ws.Range("E21").Value = '  +2.00%  '

$ws.Range("D22").Value = '14.49'
$ws.Range("E22").Value = '  +1.62%  '

$ws.Range("D23").Value = '11.15'
$ws.Range("E23").Value = '  +3.37%  '

$ws.Range("D24").Value = '2.282'
$ws.Range("E24").Value = '  +1.53%  '

$ws.Range("D25").Value = '20.588.11'
$ws.Range("E25").Value = '  +1.55%  '

$ws.Range("D26").Value = '147.11'
$ws.Range("E26").Value = '  +5.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.080'
$ws.Range("E27").Value = '  +2.01%  '

$ws.Range("D28").Value = '17.16'
$ws.Range("E28").Value = '  +0.96%  '

$ws.Range("D29").Value = '1.632.54'
$ws.Range("E29").Value = '  +2.59%  '

$ws.Range("D30").Value = '114.52'
$ws.Range("E30").Value = '  +3.40%  '

$ws.Range("D31").Value = '3.948'
$ws.Range("E31").Value = '  -1.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.930'
$ws.Range("E32").Value = '  +2.28%  '

$ws.Range("D33").Value = '0.07913'
$ws.Range("E33").Value = '  +2.88%  '

$ws.Range("D34").Value = '0.7963'
$ws.Range("E34").Value = '  +1.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.200'
$ws.Range("E35").Value = '  +7.48%  '

$ws.Range("D36").Value = '1.437'
$ws.Range("E36").Value = '  -1.18%  '

$ws.Range("D37").Value = '0.05681'
$ws.Range("E37").Value = '  +0.27%  '

$ws.Range("D38").Value = '4.697'
$ws.Range("E38").Value = '  +0.91%  '

$ws.Range("E39").Value = '  +2.28%  '

$ws.Range("D40").Value = '0.02018'
$ws.Range("E40").Value = '  +1.18%  '

$ws.Range("D41").Value = '10.26'
$ws.Range("E41").Value = '  +0.72%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '7.391'
$ws.Range("E42").Value = '  +5.71%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '0.1845'
$ws.Range("E43").Value = '  +0.20%  '

$ws.Range("D44").Value = '3.515'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5220'
$ws.Range("E45").Value = '  +0.25%  '

$ws.Range("D46").Value = '11.93'
$ws.Range("E46").Value = '  +0.85%  '

$ws.Range("D47").Value = '120.06'
$ws.Range("E47").Value = '  +3.95%  '

$ws.Range("D48").Value = '0.5158'
$ws.Range("E48").Value = '  +1.21%  '

$ws.Range("D49").Value = '1.809'
$ws.Range("E49").Value = '  +4.51%  '

$ws.Range("D50").Value = '0.06404'
$ws.Range("E50").Value = '  +0.57%  '

$ws.Range("D51").Value = '0.9942'
$ws.Range("E51").Value = '  +0.57%  '
